$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.548.80"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").Value = "3.230.88"
$ws.Range("E3").Value = "  +1.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.15"
$ws.Range("E5").Value = "  +1.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.36"
$ws.Range("E6").Value = "  +3.48%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.230.59"
$ws.Range("E8").Value = "  +1.69%  "

$ws.Range("E9").Value = "  +2.42%  "

$ws.Range("E10").Value = "  +1.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.68"
$ws.Range("E11").Value = "  -6.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.509"
$ws.Range("E12").Value = "  -0.71%  "

$ws.Range("E13").Value = "  +2.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.16"
$ws.Range("E14").Value = "  +0.90%  "

$ws.Range("D15").Value = "3.759.74"
$ws.Range("E15").Value = "  +1.74%  "

$ws.Range("D16").Value = "66.665.42"
$ws.Range("E16").Value = "  +1.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.48"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").Value = "3.237.54"
$ws.Range("E18").Value = "  +1.94%  "

$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.38"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.32"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.09"
$ws.Range("E23").Value = "  +1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.80"
$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.84"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.18"
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("E29").Value = "  +5.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.01"
$ws.Range("E30").Value = "  +4.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.07"
$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.25"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.18"
$ws.Range("E34").Value = "  -3.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.53"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "516.99"
$ws.Range("E36").Value = "  +7.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0959"
$ws.Range("E37").Value = "  +6.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.25"
$ws.Range("E38").Value = "  +2.79%  "

$ws.Range("D39").Value = "0.0₃0773"
$ws.Range("E39").Value = "  +18.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0420"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.03"
$ws.Range("E41").Value = "  +6.74%  "

$ws.Range("E42").Value = "  +6.14%  "

$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.300"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  +2.73%  "

$ws.Range("D46").Value = "2.879.85"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.63"
$ws.Range("E47").Value = "  +1.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  +4.84%  "

$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("E51").Value = "  +2.45%  "
